# Apply scheduled market-price/profit updates to the Ixion leve-profit tracking workbook.
# For each affected row, update currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (H/I/J), LevePriceNQ / LevePriceHQ (K/L) and the resulting LeveProfitNQ / LeveProfitHQ (M/N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 799.4
$ws.Range("I19").Value = 400
$ws.Range("K19").Value = 400
$ws.Range("M19").Value = -225

# Row 106
$ws.Range("H106").Value = 111114370
$ws.Range("I106").Value = 55559576
$ws.Range("J106").Value = 166669170
$ws.Range("K106").Value = 55559576
$ws.Range("L106").Value = 166669170
$ws.Range("M106").Value = -55558945
$ws.Range("N106").Value = -166670432

# Row 111
$ws.Range("H111").Value = 93245.55
$ws.Range("I111").Value = 2500
$ws.Range("J111").Value = 252050.25
$ws.Range("K111").Value = 7500
$ws.Range("L111").Value = 756150.75
$ws.Range("M111").Value = -4433
$ws.Range("N111").Value = -762284.75

# Row 118
$ws.Range("H118").Value = 950.7778
$ws.Range("I118").Value = 759.5
$ws.Range("J118").Value = 1333.3334
$ws.Range("K118").Value = 2278.5
$ws.Range("L118").Value = 4000.0002
$ws.Range("M118").Value = -621.5
$ws.Range("N118").Value = -7314.0002

# Row 127
$ws.Range("H127").Value = 802
$ws.Range("I127").Value = 416.69232
$ws.Range("J127").Value = 2054.25
$ws.Range("K127").Value = 1250.07696
$ws.Range("L127").Value = 6162.75
$ws.Range("M127").Value = 3709.92304
$ws.Range("N127").Value = -16082.75

# Row 137
$ws.Range("H137").Value = 1094.8134
$ws.Range("I137").Value = 908.4918
$ws.Range("J137").Value = 1906.6428
$ws.Range("K137").Value = 2725.4754
$ws.Range("L137").Value = 5719.928400000001
$ws.Range("M137").Value = -175.4754000000003
$ws.Range("N137").Value = -10819.9284

# Row 138
$ws.Range("H138").Value = 2669.6611
$ws.Range("I138").Value = 1089.6154
$ws.Range("J138").Value = 5750.75
$ws.Range("K138").Value = 3268.8462
$ws.Range("L138").Value = 17252.25
$ws.Range("M138").Value = 1871.1538
$ws.Range("N138").Value = -27532.25

# Row 141
$ws.Range("H141").Value = 1498.8043
$ws.Range("I141").Value = 1130.1666
$ws.Range("J141").Value = 2825.9
$ws.Range("K141").Value = 3390.4998
$ws.Range("L141").Value = 8477.700000000001
$ws.Range("M141").Value = 1789.5002
$ws.Range("N141").Value = -18837.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1503.77
$ws.Range("I32").Value = 1402.5269
$ws.Range("J32").Value = 2848.8572
$ws.Range("K32").Value = 1402.5269
$ws.Range("L32").Value = 2848.8572
$ws.Range("M32").Value = -1115.5269
$ws.Range("N32").Value = -3422.8572

# Row 61
$ws.Range("H61").Value = 3127.9019
$ws.Range("I61").Value = 3314.6135
$ws.Range("J61").Value = 1954.2858
$ws.Range("K61").Value = 3314.6135
$ws.Range("L61").Value = 1954.2858
$ws.Range("M61").Value = -3102.6135
$ws.Range("N61").Value = -2378.2858

# Row 74
$ws.Range("H74").Value = 1151.1316
$ws.Range("I74").Value = 1189.1364
$ws.Range("J74").Value = 1098.875
$ws.Range("K74").Value = 1189.1364
$ws.Range("L74").Value = 1098.875
$ws.Range("M74").Value = -315.1364000000001
$ws.Range("N74").Value = -2846.875

# Row 77
$ws.Range("H77").Value = 1151.1316
$ws.Range("I77").Value = 1189.1364
$ws.Range("J77").Value = 1098.875
$ws.Range("K77").Value = 5945.682000000001
$ws.Range("L77").Value = 5494.375
$ws.Range("M77").Value = -1577.682000000001
$ws.Range("N77").Value = -14230.375

# Row 132
$ws.Range("H132").Value = 2275412.8
$ws.Range("I132").Value = 1998.5758
$ws.Range("J132").Value = 9095655
$ws.Range("K132").Value = 5995.7274
$ws.Range("L132").Value = 27286965
$ws.Range("M132").Value = -3465.7274
$ws.Range("N132").Value = -27292025

# Row 136
$ws.Range("H136").Value = 3127.9019
$ws.Range("I136").Value = 3314.6135
$ws.Range("J136").Value = 1954.2858
$ws.Range("K136").Value = 9943.8405
$ws.Range("L136").Value = 5862.857400000001
$ws.Range("M136").Value = -7393.8405
$ws.Range("N136").Value = -10962.8574

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2729.5
$ws.Range("I134").Value = 2779.3076
$ws.Range("J134").Value = 2513.6667
$ws.Range("K134").Value = 8337.9228
$ws.Range("L134").Value = 7541.000100000001
$ws.Range("M134").Value = -5802.9228
$ws.Range("N134").Value = -12611.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 424.08334
$ws.Range("I22").Value = 453
$ws.Range("K22").Value = 453
$ws.Range("M22").Value = -103

# Row 31
$ws.Range("H31").Value = 5690.338
$ws.Range("I31").Value = 1526.7059
$ws.Range("J31").Value = 13857.462
$ws.Range("K31").Value = 1526.7059
$ws.Range("L31").Value = 13857.462
$ws.Range("M31").Value = -1231.7059
$ws.Range("N31").Value = -14447.462

# Row 34
$ws.Range("H34").Value = 5690.338
$ws.Range("I34").Value = 1526.7059
$ws.Range("J34").Value = 13857.462
$ws.Range("K34").Value = 1526.7059
$ws.Range("L34").Value = 13857.462
$ws.Range("M34").Value = -1324.7059
$ws.Range("N34").Value = -14261.462

# Row 58
$ws.Range("H58").Value = 911.7162
$ws.Range("I58").Value = 554.63464
$ws.Range("J58").Value = 1755.7273
$ws.Range("K58").Value = 554.63464
$ws.Range("L58").Value = 1755.7273
$ws.Range("M58").Value = -351.63464
$ws.Range("N58").Value = -2161.7273

# Row 107
$ws.Range("H107").Value = 20834184
$ws.Range("I107").Value = 27778330
$ws.Range("J107").Value = 1748.75
$ws.Range("K107").Value = 27778330
$ws.Range("L107").Value = 1748.75
$ws.Range("M107").Value = -27776410
$ws.Range("N107").Value = -5588.75

# Row 132
$ws.Range("H132").Value = 1766.8431
$ws.Range("I132").Value = 1397.921
$ws.Range("K132").Value = 4193.763
$ws.Range("M132").Value = -1663.763

# Row 134
$ws.Range("H134").Value = 1643.0886
$ws.Range("I134").Value = 1993.3043
$ws.Range("J134").Value = 1154.909
$ws.Range("K134").Value = 5979.9129
$ws.Range("L134").Value = 3464.727
$ws.Range("M134").Value = -3444.9129
$ws.Range("N134").Value = -8534.727000000001

# Row 136
$ws.Range("H136").Value = 911.7162
$ws.Range("I136").Value = 554.63464
$ws.Range("J136").Value = 1755.7273
$ws.Range("K136").Value = 1663.90392
$ws.Range("L136").Value = 5267.1819
$ws.Range("M136").Value = 886.09608
$ws.Range("N136").Value = -10367.1819

# Row 141
$ws.Range("H141").Value = 33756.668
$ws.Range("J141").Value = 33756.668
$ws.Range("L141").Value = 33756.668
$ws.Range("N141").Value = -44116.668

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 41818.184
$ws.Range("J37").Value = 41818.184
$ws.Range("L37").Value = 125454.552
$ws.Range("N37").Value = -125678.552

# Row 122
$ws.Range("H122").Value = 3839.3794
$ws.Range("I122").Value = 420.6875
$ws.Range("J122").Value = 8047
$ws.Range("K122").Value = 3786.1875
$ws.Range("L122").Value = 72423
$ws.Range("M122").Value = -1336.1875
$ws.Range("N122").Value = -77323

# Row 129
$ws.Range("H129").Value = 18519544
$ws.Range("I129").Value = 33333936
$ws.Range("J129").Value = 1556.125
$ws.Range("K129").Value = 100001808
$ws.Range("L129").Value = 4668.375
$ws.Range("M129").Value = -99996808
$ws.Range("N129").Value = -14668.375

# Row 139
$ws.Range("H139").Value = 4748
$ws.Range("I139").Value = 7183.75
$ws.Range("J139").Value = 2892.1904
$ws.Range("K139").Value = 21551.25
$ws.Range("L139").Value = 8676.5712
$ws.Range("M139").Value = -16411.25
$ws.Range("N139").Value = -18956.5712

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

# Row 122
$ws.Range("H122").Value = 2819910.2
$ws.Range("I122").Value = 4052634.5
$ws.Range("J122").Value = 2254.2856
$ws.Range("K122").Value = 12157903.5
$ws.Range("L122").Value = 6762.8568
$ws.Range("M122").Value = -12155453.5
$ws.Range("N122").Value = -11662.8568

# Row 132
$ws.Range("H132").Value = 1817.841
$ws.Range("I132").Value = 1346.7188
$ws.Range("J132").Value = 3074.1667
$ws.Range("K132").Value = 4040.1564
$ws.Range("L132").Value = 9222.500100000001
$ws.Range("M132").Value = -1510.1564
$ws.Range("N132").Value = -14282.5001

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 9898495
$ws.Range("I132").Value = 14065436
$ws.Range("J132").Value = 2012.4375
$ws.Range("K132").Value = 42196308
$ws.Range("L132").Value = 6037.3125
$ws.Range("M132").Value = -42193778
$ws.Range("N132").Value = -11097.3125

# Row 133
$ws.Range("H133").Value = 40326
$ws.Range("J133").Value = 40326
$ws.Range("L133").Value = 40326
$ws.Range("N133").Value = -45386

# Row 136
$ws.Range("H136").Value = 5278.7
$ws.Range("I136").Value = 3255.1724
$ws.Range("J136").Value = 15059.083
$ws.Range("K136").Value = 9765.5172
$ws.Range("L136").Value = 45177.249
$ws.Range("M136").Value = -7215.5172
$ws.Range("N136").Value = -50277.249

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 37500
$ws.Range("J92").Value = 37500
$ws.Range("L92").Value = 37500
$ws.Range("N92").Value = -42492

# Row 126
$ws.Range("H126").Value = 842.2857
$ws.Range("I126").Value = 459.2
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 1377.6
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = 1092.4
$ws.Range("N126").Value = -10340

# Row 132
$ws.Range("H132").Value = 857.9231
$ws.Range("I132").Value = 576.7436
$ws.Range("J132").Value = 1701.4615
$ws.Range("K132").Value = 1730.2308
$ws.Range("L132").Value = 5104.3845
$ws.Range("M132").Value = 799.7692
$ws.Range("N132").Value = -10164.3845

# Row 136
$ws.Range("H136").Value = 8066861.5
$ws.Range("I136").Value = 2449.9092
$ws.Range("J136").Value = 27779866
$ws.Range("K136").Value = 7349.7276
$ws.Range("L136").Value = 83339598
$ws.Range("M136").Value = -4799.7276
$ws.Range("N136").Value = -83344698
